$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.989.77'
$ws.Range('E2').Value = '  -0.91%  '
$ws.Range('D3').Value = '2.648.60'
$ws.Range('E3').Value = '  +0.29%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '581.39'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -0.25%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '156.52'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -0.63%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.622'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -3.39%  '
$ws.Range('D9').Value = '2.647.59'
$ws.Range('E9').Value = '  +0.32%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.119'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -3.38%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.82'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -0.01%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.385'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -1.66%  '
$ws.Range('E13').Value = '  +0.85%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '28.60'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -0.42%  '
$ws.Range('D15').Value = '3.122.21'
$ws.Range('E15').Value = '  +0.04%  '
$ws.Range('E16').Value = '  -1.43%  '
$ws.Range('D17').Value = '63.863.47'
$ws.Range('E17').Value = '  -0.76%  '
$ws.Range('D18').Value = '2.655.09'
$ws.Range('E18').Value = '  +0.54%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.21'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -0.49%  '
$ws.Range('E20').Value = '  +3.42%  '
$ws.Range('E21').Value = '  -3.41%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '346.72'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -0.42%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '68.22'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -0.03%  '
$ws.Range('E25').Value = '  +4.43%  '
$ws.Range('E26').Value = '  +0.27%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.30'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -1.46%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '587.14'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -1.70%  '
$ws.Range('E29').Value = '  +0.93%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.26'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +2.71%  '
$ws.Range('E31').Value = '  +0.08%  '
$ws.Range('E32').Value = '  -0.12%  '
$ws.Range('E33').Value = '  -0.70%  '
$ws.Range('E34').Value = '  -0.28%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '6.63'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -1.20%  '
$ws.Range('E36').Value = '  +3.25%  '
$ws.Range('E37').Value = '  -2.54%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '19.79'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -1.38%  '
$ws.Range('E39').Value = '  -0.02%  '
$ws.Range('E40').Value = '  -0.23%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '151.25'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -2.19%  '
$ws.Range('B42').Value = 'dogwifhat'
$ws.Range('C42').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.58'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +5.85%  '
$ws.Range('B43').Value = 'USDe'
$ws.Range('C43').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.999'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -0.05%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '41.92'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -0.87%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '163.13'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +2.77%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '24.45'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +3.78%  '
$ws.Range('E47').Value = '  -2.68%  '
$ws.Range('E48').Value = '  -2.25%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.637'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -0.01%  '
$ws.Range('E50').Value = '  -2.13%  '
$ws.Range('E51').Value = '  -2.08%  '
